$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6466.8823
$ws.Range("J43").Value = 7187.125
$ws.Range("L43").Value = 7187.125
$ws.Range("N43").Value = -7325.125
$ws.Range("H76").Value = 6834.375
$ws.Range("I76").Value = 4675
$ws.Range("J76").Value = 7142.857
$ws.Range("K76").Value = 4675
$ws.Range("L76").Value = 7142.857
$ws.Range("M76").Value = -4360
$ws.Range("N76").Value = -7772.857
$ws.Range("H79").Value = 6834.375
$ws.Range("I79").Value = 4675
$ws.Range("J79").Value = 7142.857
$ws.Range("K79").Value = 4675
$ws.Range("L79").Value = 7142.857
$ws.Range("M79").Value = -3583
$ws.Range("N79").Value = -9326.857
$ws.Range("H112").Value = 1563.375
$ws.Range("I112").Value = 1084.5
$ws.Range("J112").Value = 1723
$ws.Range("K112").Value = 3253.5
$ws.Range("L112").Value = 5169
$ws.Range("M112").Value = -2145.5
$ws.Range("N112").Value = -7385
$ws.Range("H138").Value = 2203.1853
$ws.Range("I138").Value = 776.11536
$ws.Range("J138").Value = 2877.8
$ws.Range("K138").Value = 2328.34608
$ws.Range("L138").Value = 8633.400000000001
$ws.Range("M138").Value = 2811.65392
$ws.Range("N138").Value = -18913.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8625047
$ws.Range("I32").Value = 10639611
$ws.Range("K32").Value = 10639611
$ws.Range("M32").Value = -10639324
$ws.Range("H97").Value = 2107.7778
$ws.Range("I97").Value = 2279.8333
$ws.Range("J97").Value = 1763.6666
$ws.Range("K97").Value = 2279.8333
$ws.Range("L97").Value = 1763.6666
$ws.Range("M97").Value = -1783.8333
$ws.Range("N97").Value = -2755.6666
$ws.Range("H123").Value = 41462.25
$ws.Range("J123").Value = 41462.25
$ws.Range("L123").Value = 41462.25
$ws.Range("N123").Value = -51262.25
$ws.Range("H138").Value = 216663
$ws.Range("J138").Value = 216663
$ws.Range("L138").Value = 216663
$ws.Range("N138").Value = -226943

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1075.25
$ws.Range("I22").Value = 1075.25
$ws.Range("K22").Value = 1075.25
$ws.Range("M22").Value = -902.25
$ws.Range("H107").Value = 3320.3333
$ws.Range("I107").Value = 4005.5
$ws.Range("K107").Value = 4005.5
$ws.Range("M107").Value = -2085.5
$ws.Range("H134").Value = 373428.53
$ws.Range("I134").Value = 3088.2727
$ws.Range("K134").Value = 9264.8181
$ws.Range("M134").Value = -6729.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 712075.5
$ws.Range("I31").Value = 3288.75
$ws.Range("J31").Value = 938887.25
$ws.Range("K31").Value = 3288.75
$ws.Range("L31").Value = 938887.25
$ws.Range("M31").Value = -2993.75
$ws.Range("N31").Value = -939477.25
$ws.Range("H34").Value = 712075.5
$ws.Range("I34").Value = 3288.75
$ws.Range("J34").Value = 938887.25
$ws.Range("K34").Value = 3288.75
$ws.Range("L34").Value = 938887.25
$ws.Range("M34").Value = -3086.75
$ws.Range("N34").Value = -939291.25
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H69").Value = 83339.8
$ws.Range("I69").Value = 20000
$ws.Range("J69").Value = 99174.75
$ws.Range("K69").Value = 20000
$ws.Range("L69").Value = 99174.75
$ws.Range("M69").Value = -19251
$ws.Range("N69").Value = -100672.75
$ws.Range("H72").Value = 83339.8
$ws.Range("I72").Value = 20000
$ws.Range("J72").Value = 99174.75
$ws.Range("K72").Value = 60000
$ws.Range("L72").Value = 297524.25
$ws.Range("M72").Value = -56256
$ws.Range("N72").Value = -305012.25
$ws.Range("H88").Value = 20222
$ws.Range("J88").Value = 20222
$ws.Range("L88").Value = 20222
$ws.Range("N88").Value = -21034
$ws.Range("H91").Value = 20222
$ws.Range("J91").Value = 20222
$ws.Range("L91").Value = 20222
$ws.Range("N91").Value = -23030
$ws.Range("H115").Value = 60994.5
$ws.Range("J115").Value = 60994.5
$ws.Range("L115").Value = 60994.5
$ws.Range("N115").Value = -63344.5
$ws.Range("H132").Value = 3779.3572
$ws.Range("I132").Value = 3325.9167
$ws.Range("K132").Value = 9977.750100000001
$ws.Range("M132").Value = -7447.750100000001
$ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2220.5454
$ws.Range("I5").Value = 2220.5454
$ws.Range("K5").Value = 6661.6362
$ws.Range("M5").Value = -6549.6362
$ws.Range("H131").Value = 23543
$ws.Range("J131").Value = 25329
$ws.Range("L131").Value = 75987
$ws.Range("N131").Value = -86067
$ws.Range("H132").Value = 2153
$ws.Range("I132").Value = 1480
$ws.Range("J132").Value = 2237.125
$ws.Range("K132").Value = 13320
$ws.Range("L132").Value = 20134.125
$ws.Range("M132").Value = -10790
$ws.Range("N132").Value = -25194.125
$ws.Range("H135").Value = 2220.5454
$ws.Range("I135").Value = 2220.5454
$ws.Range("K135").Value = 19984.9086
$ws.Range("M135").Value = -17449.9086

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15443.889
$ws.Range("I70").Value = 22599.4
$ws.Range("K70").Value = 22599.4
$ws.Range("M70").Value = -22329.4
$ws.Range("H73").Value = 15443.889
$ws.Range("I73").Value = 22599.4
$ws.Range("K73").Value = 22599.4
$ws.Range("M73").Value = -21663.4
$ws.Range("H102").Value = 2832.0667
$ws.Range("I102").Value = 2317.65
$ws.Range("J102").Value = 3860.9
$ws.Range("K102").Value = 2317.65
$ws.Range("L102").Value = 3860.9
$ws.Range("M102").Value = -695.6500000000001
$ws.Range("N102").Value = -7104.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3258.7896
$ws.Range("I46").Value = 2275.9167
$ws.Range("J46").Value = 4943.7144
$ws.Range("K46").Value = 2275.9167
$ws.Range("L46").Value = 4943.7144
$ws.Range("M46").Value = -2087.9167
$ws.Range("N46").Value = -5319.7144
$ws.Range("H93").Value = 66667750
$ws.Range("I93").Value = 125000810
$ws.Range("J93").Value = 1406.5714
$ws.Range("K93").Value = 125000810
$ws.Range("L93").Value = 1406.5714
$ws.Range("M93").Value = -124999562
$ws.Range("N93").Value = -3902.5714
$ws.Range("H132").Value = 850980.25
$ws.Range("I132").Value = 31209.666
$ws.Range("K132").Value = 93628.99800000001
$ws.Range("M132").Value = -91098.99800000001
$ws.Range("H136").Value = 49605.816
$ws.Range("I136").Value = 6556.0454
$ws.Range("J136").Value = 135705.36
$ws.Range("K136").Value = 19668.1362
$ws.Range("L136").Value = 407116.08
$ws.Range("M136").Value = -17118.1362
$ws.Range("N136").Value = -412216.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 10000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -9587
$ws.Range("H62").Value = 28580000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 28580000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H81").Value = 5357.643
$ws.Range("I81").Value = 1546.2727
$ws.Range("J81").Value = 19332.666
$ws.Range("K81").Value = 3092.5454
$ws.Range("L81").Value = 38665.332
$ws.Range("M81").Value = -2031.5454
$ws.Range("N81").Value = -40787.332
$ws.Range("H84").Value = 5357.643
$ws.Range("I84").Value = 1546.2727
$ws.Range("J84").Value = 19332.666
$ws.Range("K84").Value = 15462.727
$ws.Range("L84").Value = 193326.66
$ws.Range("M84").Value = -10158.727
$ws.Range("N84").Value = -203934.66
$ws.Range("H118").Value = 126000
$ws.Range("J118").Value = 126000
$ws.Range("L118").Value = 126000
$ws.Range("N118").Value = -129314
$ws.Range("H122").Value = 2349.16
$ws.Range("J122").Value = 3301.6667
$ws.Range("L122").Value = 9905.000100000001
$ws.Range("N122").Value = -14805.0001
$ws.Range("H132").Value = 305936.72
$ws.Range("I132").Value = 2746.9
$ws.Range("K132").Value = 8240.700000000001
$ws.Range("M132").Value = -5710.700000000001
$ws.Range("N39").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
